$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3.117862748653088
$ws.Range("D2").Value = 3.046292779653378
$ws.Range("E2").Value = 40.45818556876116
$ws.Range("F2").Value = 18.19546754165552
$ws.Range("G2").Value = 20.10098229929034
$ws.Range("H2").Value = 11.00280798847334
$ws.Range("I2").Value = 15.21992936458322
$ws.Range("M2").Value = 57.86297625733253
$ws.Range("C3").Value = 2.985404822067304
$ws.Range("D3").Value = 2.912887345307611
$ws.Range("E3").Value = 37.74110372221218
$ws.Range("F3").Value = 18.30961979242654
$ws.Range("G3").Value = 19.94606846167068
$ws.Range("H3").Value = 11.20115526189556
$ws.Range("I3").Value = 15.41042675248976
$ws.Range("M3").Value = 54.30956251382917
$ws.Range("C4").Value = 2.900653647319033
$ws.Range("D4").Value = 2.831489493402262
$ws.Range("E4").Value = 35.97246944501774
$ws.Range("F4").Value = 18.41178756805033
$ws.Range("G4").Value = 19.91466741458554
$ws.Range("H4").Value = 11.33283489561023
$ws.Range("I4").Value = 15.54940848617844
$ws.Range("M4").Value = 52.00425333440067
$ws.Range("C5").Value = 2.865274124193924
$ws.Range("D5").Value = 2.85537015173318
$ws.Range("E5").Value = 35.22632815649312
$ws.Range("F5").Value = 18.46112432416463
$ws.Range("G5").Value = 19.91730352813438
$ws.Range("H5").Value = 11.38890164584052
$ws.Range("I5").Value = 15.61132497513842
$ws.Range("M5").Value = 51.03382023237722
$ws.Range("C6").Value = 2.859349072903596
$ws.Range("D6").Value = 2.85937870655866
$ws.Range("E6").Value = 35.10089219952231
$ws.Range("F6").Value = 18.46977103922222
$ws.Range("G6").Value = 19.91865351042544
$ws.Range("H6").Value = 11.39835453862845
$ws.Range("I6").Value = 15.62191770812028
$ws.Range("M6").Value = 50.87081096585052
$ws.Range("C7").Value = 2.900179892803174
$ws.Range("D7").Value = 2.831808645075224
$ws.Range("E7").Value = 35.96250983502568
$ws.Range("F7").Value = 18.41242224167766
$ws.Range("G7").Value = 19.91464135899985
$ws.Range("H7").Value = 11.33358139235719
$ws.Range("I7").Value = 15.55022246771621
$ws.Range("M7").Value = 51.99129114053397
$ws.Range("C8").Value = 3.07291078442573
$ws.Range("D8").Value = 3.000911491442014
$ws.Range("E8").Value = 39.54205178757576
$ws.Range("F8").Value = 18.22795469900961
$ws.Range("G8").Value = 20.03403318266545
$ws.Range("H8").Value = 11.06909385686412
$ws.Range("I8").Value = 15.28088827812534
$ws.Range("M8").Value = 56.66340529209708
$ws.Range("C9").Value = 3.383904000936513
$ws.Range("D9").Value = 3.317479618262484
$ws.Range("E9").Value = 45.77360701996262
$ws.Range("F9").Value = 18.1361192416512
$ws.Range("G9").Value = 20.79200787107321
$ws.Range("H9").Value = 10.63289873603256
$ws.Range("I9").Value = 14.93901442756839
$ws.Range("M9").Value = 64.84605521001947
$ws.Range("C10").Value = 3.594945687979702
$ws.Range("D10").Value = 3.53619567158897
$ws.Range("E10").Value = 49.88316473305844
$ws.Range("F10").Value = 18.25237481984129
$ws.Range("G10").Value = 21.6809912343399
$ws.Range("H10").Value = 10.36864042461847
$ws.Range("I10").Value = 14.81730223377305
$ws.Range("M10").Value = 70.26494688010965
$ws.Range("C11").Value = 3.687112036132788
$ws.Range("D11").Value = 3.632812661679725
$ws.Range("E11").Value = 51.65348998288665
$ws.Range("F11").Value = 18.34859557229681
$ws.Range("G11").Value = 22.15693607674933
$ws.Range("H11").Value = 10.26207007294041
$ws.Range("I11").Value = 14.79338150661488
$ws.Range("M11").Value = 72.60301978195899
$ws.Range("C12").Value = 3.721459821253334
$ws.Range("D12").Value = 3.668997573977312
$ws.Range("E12").Value = 52.30982486036469
$ws.Range("F12").Value = 18.39149078529975
$ws.Range("G12").Value = 22.34721485605951
$ws.Range("H12").Value = 10.22380744151784
$ws.Range("I12").Value = 14.78910806253225
$ws.Range("M12").Value = 73.47028781423751
$ws.Range("C13").Value = 3.714087037396975
$ws.Range("D13").Value = 3.661222169617532
$ws.Range("E13").Value = 52.1690928792763
$ws.Range("F13").Value = 18.38196215233332
$ws.Range("G13").Value = 22.30579290959627
$ws.Range("H13").Value = 10.23195269415053
$ws.Range("I13").Value = 14.78981157883351
$ws.Range("M13").Value = 73.28430884111003
$ws.Range("C14").Value = 3.689948946392116
$ws.Range("D14").Value = 3.635797595300146
$ws.Range("E14").Value = 51.7077674057125
$ws.Range("F14").Value = 18.35199400792878
$ws.Range("G14").Value = 22.17239077770057
$ws.Range("H14").Value = 10.25887944546922
$ws.Range("I14").Value = 14.79293257518735
$ws.Range("M14").Value = 72.67473241755772
$ws.Range("C15").Value = 3.675091558287231
$ws.Range("D15").Value = 3.62017236098804
$ws.Range("E15").Value = 51.4233688787271
$ws.Range("F15").Value = 18.33448471799444
$ws.Range("G15").Value = 22.09197785474082
$ws.Range("H15").Value = 10.27564952064865
$ws.Range("I15").Value = 14.79547497446275
$ws.Range("M15").Value = 72.2989958561749
$ws.Range("C16").Value = 3.588844974820323
$ws.Range("D16").Value = 3.529824601196122
$ws.Range("E16").Value = 49.76549250997891
$ws.Range("F16").Value = 18.24697898963311
$ws.Range("G16").Value = 21.65130883857881
$ws.Range("H16").Value = 10.37589059713721
$ws.Range("I16").Value = 14.81952105267768
$ws.Range("M16").Value = 70.10960521788819
$ws.Range("C17").Value = 3.534950600767039
$ws.Range("D17").Value = 3.473670084828929
$ws.Range("E17").Value = 48.72319369756288
$ws.Range("F17").Value = 18.20457160054458
$ws.Range("G17").Value = 21.3991625038668
$ws.Range("H17").Value = 10.4409700212101
$ws.Range("I17").Value = 14.84251765218784
$ws.Range("M17").Value = 68.73405709940108
$ws.Range("C18").Value = 3.503590724769633
$ws.Range("D18").Value = 3.441099777747733
$ws.Range("E18").Value = 48.11434742795232
$ws.Range("F18").Value = 18.1842480158676
$ws.Range("G18").Value = 21.26089086518329
$ws.Range("H18").Value = 10.47967761068194
$ws.Range("I18").Value = 14.85868952368887
$ws.Range("M18").Value = 67.93091273850054
$ws.Range("C19").Value = 3.492910854921754
$ws.Range("D19").Value = 3.430025024162858
$ws.Range("E19").Value = 47.90659158341467
$ws.Range("F19").Value = 18.17805785456827
$ws.Range("G19").Value = 21.21524028895621
$ws.Range("H19").Value = 10.49299856649724
$ws.Range("I19").Value = 14.86466260996481
$ws.Range("M19").Value = 67.65692247777952
$ws.Range("C20").Value = 3.540725157732099
$ws.Range("D20").Value = 3.479675868392933
$ws.Range("E20").Value = 48.83511373187852
$ws.Range("F20").Value = 18.20866313037461
$ws.Range("G20").Value = 21.42530604524976
$ws.Range("H20").Value = 10.43390925922363
$ws.Range("I20").Value = 14.83976290127931
$ws.Range("M20").Value = 68.8817238840238
$ws.Range("C21").Value = 3.697053916593846
$ws.Range("D21").Value = 3.643276211407927
$ws.Range("E21").Value = 51.84364948852161
$ws.Range("F21").Value = 18.36061950902042
$ws.Range("G21").Value = 22.21130392409385
$ws.Range("H21").Value = 10.25091250211756
$ws.Range("I21").Value = 14.79188397628014
$ws.Range("M21").Value = 72.85426992790298
$ws.Range("C22").Value = 3.795997668095413
$ws.Range("D22").Value = 3.747862931069757
$ws.Range("E22").Value = 53.72811402757507
$ws.Range("F22").Value = 18.4976500701577
$ws.Range("G22").Value = 22.78339180836969
$ws.Range("H22").Value = 10.14358213369749
$ws.Range("I22").Value = 14.78859343804165
$ws.Range("M22").Value = 75.34509059182703
$ws.Range("C23").Value = 3.743484726170014
$ws.Range("D23").Value = 3.692252353962883
$ws.Range("E23").Value = 52.72975365727019
$ws.Range("F23").Value = 18.42099973047372
$ws.Range("G23").Value = 22.47281964945533
$ws.Range("H23").Value = 10.1996979619892
$ws.Range("I23").Value = 14.78770439643253
$ws.Range("M23").Value = 74.02528606865795
$ws.Range("C24").Value = 3.538115649288525
$ws.Range("D24").Value = 3.476961546380185
$ws.Range("E24").Value = 48.78454464786194
$ws.Range("F24").Value = 18.20680073089208
$ws.Range("G24").Value = 21.41346569327565
$ws.Range("H24").Value = 10.43709741352301
$ws.Range("I24").Value = 14.84099915312579
$ws.Range("M24").Value = 68.81500212681868
$ws.Range("C25").Value = 3.302808276052957
$ws.Range("D25").Value = 3.234312358067003
$ws.Range("E25").Value = 44.17097840077794
$ws.Range("F25").Value = 18.12998083118707
$ws.Range("G25").Value = 20.52893027953693
$ws.Range("H25").Value = 10.74148834211996
$ws.Range("I25").Value = 15.0099124396435
$ws.Range("M25").Value = 62.73699630211868
